# Update of all scripts and data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Raja asterias, 1-RAP): weight bump + new Numb value
$ws.Range("G2").Value = 0.061
$ws.Range("H2").Value = 1

# Row 6 (Raja asterias, 2-RAP): weight bump + new Numb value
$ws.Range("G6").Value = 0.061
$ws.Range("H6").Value = 1

# Rows 9-38: refresh the RF (raising factor) column
for ($r = 9; $r -le 38; $r++) {
    $ws.Cells.Item($r, 9).Value = 7.100583333333333
}

# Discard-style rows (Numb was 0) flip to -1
$ws.Range("H17").Value = -1
$ws.Range("H21").Value = -1
$ws.Range("H33").Value = -1
$ws.Range("H38").Value = -1
